$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = 65

$ws.Range("A$r").Value = "2542. Maximum Subsequence Score"

$ws.Range("B$r").Value = "Medium"
$ws.Range("B$r").Interior.Color = $ws.Range("B64").Interior.Color

$ws.Range("C$r").Value = "Heaps"

$ws.Range("D$r").Value = "The crux is that we consider for each element of nums2[i] as aminimum once and check for max possible values in num1. We first sort nums2 and keep their mappings to nums1, then iterate the pairs maintaing a heap with size k. First form the pairs as int[][]. Iterate the pairs and sort with nums2 from large to small. Keep a priority queue of size k. Each time we introduce a new pair, maintain the min value of nums2 and the sum of the priority queue. If the size of the q > k, pop min nums1, update sum -= nums1. If the size of q = k, update res = max(res, sum * nums2[i])."

$linkUrl = "https://leetcode.com/problems/maximum-subsequence-score/solutions/3557445/java-solution-for-maximum-subsequence-score-problem/?envType=study-plan-v2&envId=leetcode-75"
$ws.Range("E$r").Value = "$linkUrl "
$ws.Hyperlinks.Add($ws.Range("E$r"), $linkUrl)
$ws.Range("E64").Copy()
$ws.Range("E$r").PasteSpecial(-4122)

$ws.Application.CutCopyMode = 0

$ws.Range("E68").Select() | Out-Null
